$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -6.491699999999994
$ws.Range("A3").Value = -21.19850000000003
$ws.Range("D5").Value = -8.023599999999998
$ws.Range("E5").Value = 13.15019999999999
$ws.Range("E9").Value = 14.96020000000001
$ws.Range("E11").Value = 13.4881
$ws.Range("A14").Value = -20.61199999999998
$ws.Range("A16").Value = -20.54199999999999
$ws.Range("D16").Value = -8.104599999999996
$ws.Range("E17").Value = 13.83230000000001
$ws.Range("A21").Value = -21.19639999999999
$ws.Range("E21").Value = 12.6176
$ws.Range("A23").Value = -21.41380000000002
$ws.Range("A25").Value = -22.44870000000003
